$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data rows 2-5 with new values (diff replaces all numeric cells) ---
# Row 2
$ws.Cells.Item(2, 1).Value = 45105.50694444445
$ws.Cells.Item(2, 2).Value = 20.658
$ws.Cells.Item(2, 3).Value = 14.028
$ws.Cells.Item(2, 4).Value = 4.08
$ws.Cells.Item(2, 5).Value = 43.704
$ws.Cells.Item(2, 6).Value = 35.637
$ws.Cells.Item(2, 7).Value = 16.257
$ws.Cells.Item(2, 8).Value = 52.557
$ws.Cells.Item(2, 9).Value = 25.014
$ws.Cells.Item(2, 10).Value = 10.522
$ws.Cells.Item(2, 11).Value = 16.084
$ws.Cells.Item(2, 12).Value = 17.263
$ws.Cells.Item(2, 13).Value = 18.024
$ws.Cells.Item(2, 14).Value = 5.19
$ws.Cells.Item(2, 15).Value = 16.166
$ws.Cells.Item(2, 16).Value = 22.644
$ws.Cells.Item(2, 17).Value = 13.689
$ws.Cells.Item(2, 18).Value = 3.709
$ws.Cells.Item(2, 19).Value = 2.453
$ws.Cells.Item(2, 20).Value = 238.72
$ws.Cells.Item(2, 21).Value = 44.979
$ws.Cells.Item(2, 22).Value = 14.922
$ws.Cells.Item(2, 23).Value = 29.658
$ws.Cells.Item(2, 24).Value = 15.316
$ws.Cells.Item(2, 25).Value = 2.936
$ws.Cells.Item(2, 26).Value = 26.039
$ws.Cells.Item(2, 27).Value = 13.181
$ws.Cells.Item(2, 28).Value = 11.919
$ws.Cells.Item(2, 29).Value = 13.951
$ws.Cells.Item(2, 30).Value = 17.711
$ws.Cells.Item(2, 31).Value = 3.458
$ws.Cells.Item(2, 32).Value = 46.546
$ws.Cells.Item(2, 33).Value = 8.337999999999999
$ws.Cells.Item(2, 34).Value = 18.656

# Row 3
$ws.Cells.Item(3, 1).Value = 45105.51388888889
$ws.Cells.Item(3, 2).Value = 21.619
$ws.Cells.Item(3, 3).Value = 15.538
$ws.Cells.Item(3, 4).Value = 1.983
$ws.Cells.Item(3, 5).Value = 46.625
$ws.Cells.Item(3, 6).Value = 38.224
$ws.Cells.Item(3, 7).Value = 17.013
$ws.Cells.Item(3, 8).Value = 65.31399999999999
$ws.Cells.Item(3, 9).Value = 26.178
$ws.Cells.Item(3, 10).Value = 11.432
$ws.Cells.Item(3, 11).Value = 17.079
$ws.Cells.Item(3, 12).Value = 18.703
$ws.Cells.Item(3, 13).Value = 19.594
$ws.Cells.Item(3, 14).Value = 5.435
$ws.Cells.Item(3, 15).Value = 16.918
$ws.Cells.Item(3, 16).Value = 23.946
$ws.Cells.Item(3, 17).Value = 14.39
$ws.Cells.Item(3, 18).Value = 1.647
$ws.Cells.Item(3, 19).Value = 1.234
$ws.Cells.Item(3, 20).Value = 250.204
$ws.Cells.Item(3, 21).Value = 47.326
$ws.Cells.Item(3, 22).Value = 15.616
$ws.Cells.Item(3, 23).Value = 31.557
$ws.Cells.Item(3, 24).Value = 16.606
$ws.Cells.Item(3, 25).Value = 2.723
$ws.Cells.Item(3, 26).Value = 31.616
$ws.Cells.Item(3, 27).Value = 13.794
$ws.Cells.Item(3, 28).Value = 12.374
$ws.Cells.Item(3, 29).Value = 14.515
$ws.Cells.Item(3, 30).Value = 19.399
$ws.Cells.Item(3, 31).Value = 1.266
$ws.Cells.Item(3, 32).Value = 59.3
$ws.Cells.Item(3, 33).Value = 8.77
$ws.Cells.Item(3, 34).Value = 19.524

# Row 4
$ws.Cells.Item(4, 1).Value = 45105.52083333334
$ws.Cells.Item(4, 2).Value = 21.139
$ws.Cells.Item(4, 3).Value = 15.421
$ws.Cells.Item(4, 4).Value = 1.479
$ws.Cells.Item(4, 5).Value = 45.758
$ws.Cells.Item(4, 6).Value = 37.566
$ws.Cells.Item(4, 7).Value = 16.636
$ws.Cells.Item(4, 8).Value = 65.17100000000001
$ws.Cells.Item(4, 9).Value = 25.596
$ws.Cells.Item(4, 10).Value = 11.287
$ws.Cells.Item(4, 11).Value = 16.796
$ws.Cells.Item(4, 12).Value = 18.394
$ws.Cells.Item(4, 13).Value = 19.308
$ws.Cells.Item(4, 14).Value = 5.314
$ws.Cells.Item(4, 15).Value = 16.542
$ws.Cells.Item(4, 16).Value = 23.487
$ws.Cells.Item(4, 17).Value = 14.01
$ws.Cells.Item(4, 18).Value = 1.112
$ws.Cells.Item(4, 19).Value = 0.961
$ws.Cells.Item(4, 20).Value = 244.483
$ws.Cells.Item(4, 21).Value = 46.285
$ws.Cells.Item(4, 22).Value = 15.269
$ws.Cells.Item(4, 23).Value = 30.994
$ws.Cells.Item(4, 24).Value = 16.325
$ws.Cells.Item(4, 25).Value = 2.535
$ws.Cells.Item(4, 26).Value = 31.365
$ws.Cells.Item(4, 27).Value = 13.487
$ws.Cells.Item(4, 28).Value = 12.028
$ws.Cells.Item(4, 29).Value = 14.123
$ws.Cells.Item(4, 30).Value = 19.164
$ws.Cells.Item(4, 31).Value = 0.773
$ws.Cells.Item(4, 32).Value = 59.107
$ws.Cells.Item(4, 33).Value = 8.595000000000001
$ws.Cells.Item(4, 34).Value = 19.09

# Row 5
$ws.Cells.Item(5, 1).Value = 45105.52777777778
$ws.Cells.Item(5, 2).Value = 0.96
$ws.Cells.Item(5, 3).Value = 0.38
$ws.Cells.Item(5, 4).Value = 0.59
$ws.Cells.Item(5, 5).Value = 1.94
$ws.Cells.Item(5, 6).Value = 1.37
$ws.Cells.Item(5, 7).Value = 0.77
$ws.Cells.Item(5, 8).Value = 11.86
$ws.Cells.Item(5, 9).Value = 1.16
$ws.Cells.Item(5, 10).Value = 0.57
$ws.Cells.Item(5, 11).Value = 0.47
$ws.Cells.Item(5, 12).Value = 0.82
$ws.Cells.Item(5, 13).Value = 0.74
$ws.Cells.Item(5, 14).Value = 0.29
$ws.Cells.Item(5, 15).Value = 0.75
$ws.Cells.Item(5, 16).Value = 1.21
$ws.Cells.Item(5, 17).Value = 0.88
$ws.Cells.Item(5, 18).Value = 0.7
$ws.Cells.Item(5, 19).Value = 0.26
$ws.Cells.Item(5, 20).Value = 4.32
$ws.Cells.Item(5, 21).Value = 2.78
$ws.Cells.Item(5, 22).Value = 0.6899999999999999
$ws.Cells.Item(5, 23).Value = 1.88
$ws.Cells.Item(5, 24).Value = 0.87
$ws.Cells.Item(5, 25).Value = 0.39
$ws.Cells.Item(5, 26).Value = 4.92
$ws.Cells.Item(5, 27).Value = 0.61
$ws.Cells.Item(5, 28).Value = 0.71
$ws.Cells.Item(5, 29).Value = 0.79
$ws.Cells.Item(5, 30).Value = 0.73
$ws.Cells.Item(5, 31).Value = 0.5600000000000001
$ws.Cells.Item(5, 32).Value = 11.42
$ws.Cells.Item(5, 33).Value = 0.3
$ws.Cells.Item(5, 34).Value = 0.89

# --- Delete row 6 (data shrinks from 5 data rows to 4 data rows) ---
$ws.Rows.Item(6).Delete()

# --- Adjust column widths (ColumnWidth input = target char width - 6/7 so that
#     Excel stored width attribute rounds to the exact target integer) ---
$ws.Columns.Item(2).ColumnWidth = 7.142857142857143  # B -> width 8
$ws.Columns.Item(3).ColumnWidth = 7.142857142857143  # C -> width 8
$ws.Columns.Item(7).ColumnWidth = 7.142857142857143  # G -> width 8
$ws.Columns.Item(9).ColumnWidth = 7.142857142857143  # I -> width 8
$ws.Columns.Item(10).ColumnWidth = 7.142857142857143  # J -> width 8
$ws.Columns.Item(11).ColumnWidth = 7.142857142857143  # K -> width 8
$ws.Columns.Item(12).ColumnWidth = 7.142857142857143  # L -> width 8
$ws.Columns.Item(13).ColumnWidth = 7.142857142857143  # M -> width 8
$ws.Columns.Item(15).ColumnWidth = 7.142857142857143  # O -> width 8
$ws.Columns.Item(16).ColumnWidth = 7.142857142857143  # P -> width 8
$ws.Columns.Item(17).ColumnWidth = 7.142857142857143  # Q -> width 8
$ws.Columns.Item(20).ColumnWidth = 8.142857142857142  # T -> width 9
$ws.Columns.Item(22).ColumnWidth = 7.142857142857143  # V -> width 8
$ws.Columns.Item(23).ColumnWidth = 7.142857142857143  # W -> width 8
$ws.Columns.Item(24).ColumnWidth = 7.142857142857143  # X -> width 8
$ws.Columns.Item(26).ColumnWidth = 7.142857142857143  # Z -> width 8
$ws.Columns.Item(27).ColumnWidth = 7.142857142857143  # AA -> width 8
$ws.Columns.Item(28).ColumnWidth = 7.142857142857143  # AB -> width 8
$ws.Columns.Item(29).ColumnWidth = 7.142857142857143  # AC -> width 8
$ws.Columns.Item(30).ColumnWidth = 7.142857142857143  # AD -> width 8
$ws.Columns.Item(31).ColumnWidth = 6.142857142857143  # AE -> width 7
$ws.Columns.Item(34).ColumnWidth = 7.142857142857143  # AH -> width 8
